$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# D-column values are forced to text (NumberFormat "@" then ClearFormats)
# so Excel does not silently reinterpret price strings like "1.000"/"0.5159"
# as numbers, and so no residual cell style/number-format is left behind.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '25.991.02'
$cell.ClearFormats()
$ws.Range("E2").Value = '  -0.20%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.742.94'
$cell.ClearFormats()
$ws.Range("E3").Value = '  -0.25%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.ClearFormats()
$ws.Range("E4").Value = '  +0.02%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '250.51'
$cell.ClearFormats()
$ws.Range("E5").Value = '  +7.29%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()
$ws.Range("E6").Value = '  +0.07%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.5159'
$cell.ClearFormats()
$ws.Range("E7").Value = '  -2.33%  '

$ws.Range("E8").Value = '  -0.34%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.06194'
$cell.ClearFormats()
$ws.Range("E9").Value = '  +0.20%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '1.741.84'
$cell.ClearFormats()
$ws.Range("E10").Value = '  -0.29%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07235'
$cell.ClearFormats()
$ws.Range("E11").Value = '  +0.49%  '

$ws.Range("E12").Value = '  -0.66%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.6509'
$cell.ClearFormats()
$ws.Range("E13").Value = '  +1.54%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '4.636'
$cell.ClearFormats()
$ws.Range("E14").Value = '  +0.81%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '77.81'
$cell.ClearFormats()
$ws.Range("E15").Value = '  -0.64%  '

$ws.Range("E16").Value = '  +0.02%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.ClearFormats()
$ws.Range("E17").Value = '  +0.02%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '26.013.91'
$cell.ClearFormats()
$ws.Range("E18").Value = '  +0.15%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '11.85'
$cell.ClearFormats()
$ws.Range("E19").Value = '  +2.30%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.000006810'
$cell.ClearFormats()
$ws.Range("E20").Value = '  +1.18%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '1.964.30'
$cell.ClearFormats()
$ws.Range("E21").Value = '  -0.58%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '4.289'
$cell.ClearFormats()
$ws.Range("E22").Value = '  -0.79%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '8.694'
$cell.ClearFormats()
$ws.Range("E23").Value = '  -1.25%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '5.398'
$cell.ClearFormats()
$ws.Range("E24").Value = '  +3.83%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '136.15'
$cell.ClearFormats()
$ws.Range("E25").Value = '  -2.26%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '1.513'
$cell.ClearFormats()
$ws.Range("E26").Value = '  -0.38%  '

$ws.Range("E27").Value = '  -0.05%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '1.784'
$cell.ClearFormats()
$ws.Range("E28").Value = '  -1.22%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '105.84'
$cell.ClearFormats()
$ws.Range("E29").Value = '  +1.62%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '3.965'
$cell.ClearFormats()
$ws.Range("E30").Value = '  +5.03%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.08241'
$cell.ClearFormats()
$ws.Range("E31").Value = '  -0.73%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '3.660'
$cell.ClearFormats()
$ws.Range("E32").Value = '  -0.31%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.04704'
$cell.ClearFormats()
$ws.Range("E33").Value = '  +4.02%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '2.655'
$cell.ClearFormats()
$ws.Range("E34").Value = '  +0.77%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.ClearFormats()
$ws.Range("E35").Value = '  +0.12%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.6265'
$cell.ClearFormats()
$ws.Range("E36").Value = '  -0.75%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.735'
$cell.ClearFormats()
$ws.Range("E37").Value = '  +1.06%  '

$ws.Range("E38").Value = '  +1.47%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '1.924'
$cell.ClearFormats()
$ws.Range("E39").Value = '  -0.29%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.ClearFormats()
$ws.Range("E40").Value = '  +0.08%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '100.57'
$cell.ClearFormats()
$ws.Range("E41").Value = '  +2.58%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.7643'
$cell.ClearFormats()
$ws.Range("E42").Value = '  +3.57%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.3868'
$cell.ClearFormats()
$ws.Range("E43").Value = '  -0.75%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '5.028'
$cell.ClearFormats()
$ws.Range("E44").Value = '  -0.29%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '6.369'
$cell.ClearFormats()
$ws.Range("E45").Value = '  +0.40%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.1135'
$cell.ClearFormats()
$ws.Range("E46").Value = '  -0.42%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '55.64'
$cell.ClearFormats()
$ws.Range("E47").Value = '  +3.02%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.05227'
$cell.ClearFormats()
$ws.Range("E48").Value = '  -2.18%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '30.80'
$cell.ClearFormats()
$ws.Range("E49").Value = '  +0.71%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '7.584'
$cell.ClearFormats()
$ws.Range("E50").Value = '  -1.10%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.3435'
$cell.ClearFormats()
$ws.Range("E51").Value = '  -0.66%  '
